# Updated symbol list - applies the coin-table refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is a cell reference (A1-style) and its new text value.
# Values are written as plain text (NumberFormat "@") so price/
# percentage/index strings like "244.14", "-0.78%", "11" stay text
# cells instead of being auto-coerced into numbers/percentages by Excel,
# matching the original inline-string cell type. ClearFormats() afterwards
# drops the transient text number-format so no stray style is left behind.
$updates = @(
    @{ Cell = 'D2'; Value = '244.14' }
    @{ Cell = 'E2'; Value = '-0.78%' }
    @{ Cell = 'G2'; Value = '11' }
    @{ Cell = 'D3'; Value = '27.15' }
    @{ Cell = 'E3'; Value = '3.93%' }
    @{ Cell = 'G3'; Value = '11' }
    @{ Cell = 'E4'; Value = '1.26%' }
    @{ Cell = 'G4'; Value = '11' }
    @{ Cell = 'D5'; Value = '0.05631' }
    @{ Cell = 'E5'; Value = '0.38%' }
    @{ Cell = 'G5'; Value = '11' }
    @{ Cell = 'E6'; Value = '0.02%' }
    @{ Cell = 'G6'; Value = '11' }
    @{ Cell = 'D7'; Value = '0.8162' }
    @{ Cell = 'E7'; Value = '0.55%' }
    @{ Cell = 'G7'; Value = '11' }
    @{ Cell = 'D8'; Value = '0.8307' }
    @{ Cell = 'E8'; Value = '-1.97%' }
    @{ Cell = 'G8'; Value = '11' }
    @{ Cell = 'D9'; Value = '0.1329' }
    @{ Cell = 'E9'; Value = '-1.02%' }
    @{ Cell = 'G9'; Value = '11' }
    @{ Cell = 'D10'; Value = '0.06914' }
    @{ Cell = 'E10'; Value = '-0.88%' }
    @{ Cell = 'G10'; Value = '11' }
    @{ Cell = 'D11'; Value = '0.02942' }
    @{ Cell = 'E11'; Value = '6.06%' }
    @{ Cell = 'G11'; Value = '11' }
    @{ Cell = 'D12'; Value = '0.09399' }
    @{ Cell = 'E12'; Value = '-0.08%' }
    @{ Cell = 'G12'; Value = '11' }
    @{ Cell = 'D13'; Value = '0.001512' }
    @{ Cell = 'E13'; Value = '-0.04%' }
    @{ Cell = 'G13'; Value = '11' }
    @{ Cell = 'D14'; Value = '0.04237' }
    @{ Cell = 'E14'; Value = '-9.81%' }
    @{ Cell = 'G14'; Value = '11' }
    @{ Cell = 'D15'; Value = '0.0005996' }
    @{ Cell = 'E15'; Value = '-93.88%' }
    @{ Cell = 'G15'; Value = '11' }
    @{ Cell = 'D16'; Value = '0.006108' }
    @{ Cell = 'E16'; Value = '-0.46%' }
    @{ Cell = 'G16'; Value = '11' }
    @{ Cell = 'D17'; Value = '3.566' }
    @{ Cell = 'E17'; Value = '0.22%' }
    @{ Cell = 'G17'; Value = '11' }
    @{ Cell = 'D18'; Value = '3.018' }
    @{ Cell = 'E18'; Value = '-0.08%' }
    @{ Cell = 'G18'; Value = '11' }
    @{ Cell = 'D19'; Value = '2.309' }
    @{ Cell = 'E19'; Value = '9.02%' }
    @{ Cell = 'G19'; Value = '11' }
    @{ Cell = 'G20'; Value = '11' }
    @{ Cell = 'E21'; Value = '-3.70%' }
    @{ Cell = 'G21'; Value = '11' }
    @{ Cell = 'E22'; Value = '-2.14%' }
    @{ Cell = 'G22'; Value = '11' }
    @{ Cell = 'D23'; Value = '3.738' }
    @{ Cell = 'E23'; Value = '-0.32%' }
    @{ Cell = 'G23'; Value = '11' }
    @{ Cell = 'D24'; Value = '0.1373' }
    @{ Cell = 'G24'; Value = '11' }
    @{ Cell = 'E25'; Value = '-1.71%' }
    @{ Cell = 'G25'; Value = '11' }
    @{ Cell = 'D26'; Value = '0.004487' }
    @{ Cell = 'E26'; Value = '-2.87%' }
    @{ Cell = 'G26'; Value = '11' }
    @{ Cell = 'E27'; Value = '2.06%' }
    @{ Cell = 'G27'; Value = '11' }
    @{ Cell = 'E28'; Value = '-0.50%' }
    @{ Cell = 'G28'; Value = '11' }
    @{ Cell = 'G29'; Value = '11' }
    @{ Cell = 'G30'; Value = '11' }
    @{ Cell = 'G31'; Value = '11' }
    @{ Cell = 'G32'; Value = '11' }
    @{ Cell = 'G33'; Value = '11' }
    @{ Cell = 'G34'; Value = '11' }
    @{ Cell = 'G35'; Value = '11' }
    @{ Cell = 'G36'; Value = '11' }
    @{ Cell = 'G37'; Value = '11' }
    @{ Cell = 'G38'; Value = '11' }
    @{ Cell = 'G39'; Value = '11' }
    @{ Cell = 'D40'; Value = '0.03649' }
    @{ Cell = 'E40'; Value = '-0.25%' }
    @{ Cell = 'G40'; Value = '11' }
    @{ Cell = 'B41'; Value = 'KickToken' }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick' }
    @{ Cell = 'D41'; Value = '0.006047' }
    @{ Cell = 'E41'; Value = '-1.01%' }
    @{ Cell = 'G41'; Value = '11' }
    @{ Cell = 'B42'; Value = 'BKEXToken' }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk' }
    @{ Cell = 'D42'; Value = '0.1054' }
    @{ Cell = 'E42'; Value = '-0.16%' }
    @{ Cell = 'G42'; Value = '11' }
    @{ Cell = 'B43'; Value = 'CEJI' }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji' }
    @{ Cell = 'D43'; Value = '0.001799' }
    @{ Cell = 'E43'; Value = '-10.03%' }
    @{ Cell = 'G43'; Value = '11' }
    @{ Cell = 'D44'; Value = '0.008097' }
    @{ Cell = 'E44'; Value = '-6.03%' }
    @{ Cell = 'G44'; Value = '11' }
    @{ Cell = 'D45'; Value = '0.00005401' }
    @{ Cell = 'E45'; Value = '2.02%' }
    @{ Cell = 'G45'; Value = '11' }
    @{ Cell = 'E46'; Value = '-0.03%' }
    @{ Cell = 'G46'; Value = '11' }
    @{ Cell = 'D47'; Value = '0.1089' }
    @{ Cell = 'E47'; Value = '-18.09%' }
    @{ Cell = 'G47'; Value = '11' }
    @{ Cell = 'D48'; Value = '0.002639' }
    @{ Cell = 'E48'; Value = '28.80%' }
    @{ Cell = 'G48'; Value = '11' }
    @{ Cell = 'E49'; Value = '-0.03%' }
    @{ Cell = 'G49'; Value = '11' }
    @{ Cell = 'E50'; Value = '-0.03%' }
    @{ Cell = 'G50'; Value = '11' }
    @{ Cell = 'G51'; Value = '11' }
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    $c.NumberFormat = "@"
    $c.Value = $u.Value
    $c.ClearFormats()
}

